# Commit #5: cash & deposit done
# Adds bank/deposit metadata columns (property_category, category, date,
# legislator_name, legislator_id, source_file, index) to the "存款" (deposit)
# sheet, and turns row 1 into a proper field-name header row, matching the
# pattern already used on the other sheets of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

$lastRow = $ws.UsedRange.Rows.Count

# ---- Row 1: fix header row (currently a stray duplicate of the first data
# row) so that it contains field names like every other sheet ----
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"

$headers = @("property_category", "category", "date", "legislator_name", "legislator_id", "source_file", "index")
$startCol = 7   # column G
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $startCol + $i
    $cell = $ws.Cells.Item(1, $col)
    # copy formatting (bold + border) from an existing header cell, then set text
    $ws.Range("B1").Copy($cell)
    $cell.Value = $headers[$i]
}

# ---- Data rows: append the constant metadata columns to every record row ----
for ($r = 2; $r -le $lastRow; $r++) {
    $idx = $ws.Cells.Item($r, 1).Value()

    $propertyCategoryCell = $ws.Cells.Item($r, 7)
    $ws.Cells.Item($r, 2).Copy($propertyCategoryCell)
    $propertyCategoryCell.Value = "deposit"

    $categoryCell = $ws.Cells.Item($r, 8)
    $ws.Cells.Item($r, 2).Copy($categoryCell)
    $categoryCell.Value = "normal"

    # Force the date column to stay plain text instead of being auto-converted
    # to a date serial number by Excel's input parser.
    $dateCell = $ws.Cells.Item($r, 9)
    $ws.Cells.Item($r, 2).Copy($dateCell)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2012-04-20"

    $legislatorNameCell = $ws.Cells.Item($r, 10)
    $ws.Cells.Item($r, 2).Copy($legislatorNameCell)
    $legislatorNameCell.Value = "費鴻泰"

    $legislatorIdCell = $ws.Cells.Item($r, 11)
    $ws.Cells.Item($r, 2).Copy($legislatorIdCell)
    $legislatorIdCell.Value = 1365

    $sourceFileCell = $ws.Cells.Item($r, 12)
    $ws.Cells.Item($r, 2).Copy($sourceFileCell)
    $sourceFileCell.Value = "tmpe52e1"

    $indexCell = $ws.Cells.Item($r, 13)
    $ws.Cells.Item($r, 2).Copy($indexCell)
    $indexCell.Value = $idx
}
